$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo value (C2) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 555.1347118607418

# --- Sheet "Solucion": update Pedido (A) / Salida (B) columns (rows 2-61) ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$solucionData = @(
    @{Row=2; A="Pedido_1"; B="S001"},
    @{Row=3; A="Pedido_50"; B="S021"},
    @{Row=4; A="Pedido_45"; B="S041"},
    @{Row=5; A="Pedido_4"; B="S031"},
    @{Row=6; A="Pedido_9"; B="S011"},
    @{Row=7; A="Pedido_24"; B="S051"},
    @{Row=8; A="Pedido_57"; B="S042"},
    @{Row=9; A="Pedido_22"; B="S002"},
    @{Row=10; A="Pedido_10"; B="S022"},
    @{Row=11; A="Pedido_54"; B="S052"},
    @{Row=12; A="Pedido_40"; B="S012"},
    @{Row=13; A="Pedido_29"; B="S032"},
    @{Row=14; A="Pedido_6"; B="S003"},
    @{Row=15; A="Pedido_52"; B="S043"},
    @{Row=16; A="Pedido_38"; B="S023"},
    @{Row=17; A="Pedido_25"; B="S013"},
    @{Row=18; A="Pedido_27"; B="S053"},
    @{Row=19; A="Pedido_21"; B="S033"},
    @{Row=20; A="Pedido_60"; B="S004"},
    @{Row=21; A="Pedido_20"; B="S024"},
    @{Row=22; A="Pedido_5"; B="S044"},
    @{Row=23; A="Pedido_56"; B="S014"},
    @{Row=24; A="Pedido_7"; B="S054"},
    @{Row=25; A="Pedido_16"; B="S005"},
    @{Row=26; A="Pedido_36"; B="S034"},
    @{Row=27; A="Pedido_46"; B="S025"},
    @{Row=28; A="Pedido_8"; B="S045"},
    @{Row=29; A="Pedido_23"; B="S035"},
    @{Row=30; A="Pedido_2"; B="S015"},
    @{Row=31; A="Pedido_33"; B="S055"},
    @{Row=32; A="Pedido_43"; B="S026"},
    @{Row=33; A="Pedido_31"; B="S006"},
    @{Row=34; A="Pedido_41"; B="S036"},
    @{Row=35; A="Pedido_51"; B="S046"},
    @{Row=36; A="Pedido_35"; B="S016"},
    @{Row=37; A="Pedido_53"; B="S056"},
    @{Row=38; A="Pedido_17"; B="S047"},
    @{Row=39; A="Pedido_37"; B="S027"},
    @{Row=40; A="Pedido_15"; B="S007"},
    @{Row=41; A="Pedido_30"; B="S057"},
    @{Row=42; A="Pedido_48"; B="S017"},
    @{Row=43; A="Pedido_18"; B="S037"},
    @{Row=44; A="Pedido_44"; B="S048"},
    @{Row=45; A="Pedido_11"; B="S028"},
    @{Row=46; A="Pedido_34"; B="S008"},
    @{Row=47; A="Pedido_19"; B="S058"},
    @{Row=48; A="Pedido_26"; B="S018"},
    @{Row=49; A="Pedido_59"; B="S038"},
    @{Row=50; A="Pedido_3"; B="S049"},
    @{Row=51; A="Pedido_49"; B="S029"},
    @{Row=52; A="Pedido_12"; B="S009"},
    @{Row=53; A="Pedido_39"; B="S019"},
    @{Row=54; A="Pedido_32"; B="S059"},
    @{Row=55; A="Pedido_42"; B="S039"},
    @{Row=56; A="Pedido_28"; B="S010"},
    @{Row=57; A="Pedido_13"; B="S050"},
    @{Row=58; A="Pedido_47"; B="S030"},
    @{Row=59; A="Pedido_58"; B="S020"},
    @{Row=60; A="Pedido_14"; B="S060"},
    @{Row=61; A="Pedido_55"; B="S040"}
)

foreach ($item in $solucionData) {
    $wsSolucion.Cells.Item($item.Row, 1).Value = $item.A
    $wsSolucion.Cells.Item($item.Row, 2).Value = $item.B
}

# --- Sheet "Metricas": update Tiempo values (B2:B4) ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 549.6654395069736
$wsMetricas.Range("B3").Value = 552.8009136122823
$wsMetricas.Range("B4").Value = 555.1347118607418

Write-Host "Edit complete"
